$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching style of existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-80: I column and J column values
$data = @(
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 10),
    @(9, 9),
    @(9, 10),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(9, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(6, 6),
    @(5, 5),
    @(6, 6),
    @(4, 4),
    @(5, 5),
    @(5, 5)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
